$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Cell G7 already carries the "Ruim" (red / Failed-style) cell format used throughout
# the sheet for Status="Failed" and populated "Defect ID" cells. Use it as a format
# donor via Copy + PasteSpecial(xlPasteFormats) so we reuse the existing style index
# instead of synthesizing a brand-new one.
$xlPasteFormats = -4122
$fmtSource = $ws.Range("G7")

# --- Row 4 (TC1): Status -> Failed, Defect ID -> BUG001, Notes -> new bug note ---
$fmtSource.Copy()
$ws.Range("F4").PasteSpecial($xlPasteFormats)
$ws.Range("F4").Value = "Failed"

$fmtSource.Copy()
$ws.Range("G4").PasteSpecial($xlPasteFormats)
$ws.Range("G4").Value = "BUG001"

$ws.Range("H4").Value = "Reservation's day displayed in reservations list is incorrect"
$ws.Rows.Item(4).RowHeight = 34

# --- Row 7 (TC4): Defect ID BUG001 -> BUG002 ---
$ws.Range("G7").Value = "BUG002"

# --- Row 8 (TC5): Status -> Failed, Defect ID -> BUG001 ---
$fmtSource.Copy()
$ws.Range("F8").PasteSpecial($xlPasteFormats)
$ws.Range("F8").Value = "Failed"

$fmtSource.Copy()
$ws.Range("G8").PasteSpecial($xlPasteFormats)
$ws.Range("G8").Value = "BUG001"

# --- Row 12 (TC9): Defect ID BUG002 -> BUG003 ---
$ws.Range("G12").Value = "BUG003"

# --- Row 13 (TC10): Defect ID -> BUG004 ---
$fmtSource.Copy()
$ws.Range("G13").PasteSpecial($xlPasteFormats)
$ws.Range("G13").Value = "BUG004"

# --- Row 14 (TC11): Defect ID -> BUG004 ---
$fmtSource.Copy()
$ws.Range("G14").PasteSpecial($xlPasteFormats)
$ws.Range("G14").Value = "BUG004"

# --- Row 20 (TC17): Defect ID -> BUG008 ---
$fmtSource.Copy()
$ws.Range("G20").PasteSpecial($xlPasteFormats)
$ws.Range("G20").Value = "BUG008"

# --- Row 16 (TC13): Defect ID -> BUG007 ---
$fmtSource.Copy()
$ws.Range("G16").PasteSpecial($xlPasteFormats)
$ws.Range("G16").Value = "BUG007"

# --- Row 17 (TC14): Defect ID -> BUG007 ---
$fmtSource.Copy()
$ws.Range("G17").PasteSpecial($xlPasteFormats)
$ws.Range("G17").Value = "BUG007"

$ws.Application.CutCopyMode = $false

# Update the active selection to match the document's last saved cursor position.
$ws.Activate()
$ws.Range("E14").Select()

Write-Host "edit complete"
